# Remove the "fullpathhere/" prefix from the image file name values
# stored in column A (rows 2-11 hold the shared strings that changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$rowCount = $used.Rows.Count()
for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $val -like "fullpathhere/*") {
        $cell.Value = $val.Substring(13)
    }
}

# Update the active selection on the sheet to match the saved view state.
$ws.Range("A12").Select()
